# Daily attendance processing - 2025-10-13 08:53:55
# Updates "Recorded By" email-list orderings, a handful of recomputed
# summary/coverage statistics, and flips the PHYSIOLOGY C1 Session 1 row
# (row 19) from "Pending" to "Recorded" now that 81/221 attendance has
# come in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* like a percentage (e.g. "18.2%") as
# literal text instead of letting Excel's autoconvert turn it into a
# number with a percentage format. We force the cell to Text ("@") before
# assigning, then strip the format back off (ClearFormats) and restore the
# sheet's normal centered alignment so the cell's style matches the rest
# of the table instead of staying flagged as an explicit Text format.
function Set-LiteralText {
    param($range, [string]$value)

    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4108     # xlCenter
}

# --- Recorded-by email list reorderings (same attendees, new order) ---

$ws.Range("G3").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G4").Value = "Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

$ws.Range("G12").Value = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G25").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
$ws.Range("G26").Value = "Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

$ws.Range("G34").Value = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G41").Value = "maryam.ashraf@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("G42").Value = "marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("H42").Value = "17/246"

# --- Class Statistics block (rows 6/8/9) ---

$ws.Range("L6").Value = 9
$ws.Range("L8").Value = 30
Set-LiteralText $ws.Range("L9") "20.5%"

# --- Group Statistics block (row 15/16) ---

$ws.Range("O15").Value = 4
$ws.Range("Q15").Value = 16
Set-LiteralText $ws.Range("R15") "18.2%"
Set-LiteralText $ws.Range("S15") "46.4%"

Set-LiteralText $ws.Range("S16") "30.3%"

# --- Row 19: PHYSIOLOGY C1 Session 1 recorded (was Pending) ---
# Copy the "Recorded" (green) formatting from row 3 onto row 19, then
# fill in the attendance details.

$ws.Range("A3:I3").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G19").Value = "marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("H19").Value = "81/221"
$ws.Range("I19").Value = "Recorded"
